$wb = $excel.ActiveWorkbook

$team = $wb.Worksheets.Item("Team")
$events = $wb.Worksheets.Item("Events")
$artifacts = $wb.Worksheets.Item("Artifacts")

# --- Team sheet content updates ---
$team.Range("E5").Value = "Took on Debs' work for Prish"
$team.Range("E9").Value = "Asked to work late by Cathy. Worked on Prish's item but context switched."
$team.Range("E10").Value = "Works late to test PBI 58712.
Struggled with the test environment. Had a holiday"
$team.Range("E11").Value = "Skipped one standup. Ran manual regression"
$team.Range("C15").Value = "Oli"
$team.Range("E14").Value = "Not involved"
$team.Range("E15").Value = "Not involved"
$team.Range("E16").Value = "Sometime didn't turn up to the stand-up, sometimes interrupted. Escalated issue to Cathy"
$team.Range("C17").Value = "Prish"
$team.Range("E17").Value = "Didn't attend stand-ups, and wasn't aware of scope changing. Interrupts the stand up. Talks to the team rather than through the PO"

# --- Row height adjustments (wrapped text now spans two lines) ---
$team.Rows.Item(9).RowHeight = 30
$team.Rows.Item(16).RowHeight = 30
$team.Rows.Item(17).RowHeight = 30

# --- Events sheet content updates ---
$events.Range("C5").Value = "Patty late. Some people have sporadic attendance. Cancelled once. Its OK to have lots of people at the stand up, but only the scrum team can participate, unless the team want to actively ask questions. Others should observe. Scrum team members should set daily targets and discuss how they are progressing towards the sprint goal, and what to do next."

# --- Selections / active sheet ---
$team.Activate()
$team.Range("C4").Select()
$events.Activate()
$events.Range("C3").Select()
$artifacts.Activate()
$artifacts.Range("C5").Select()
